$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 285, pre-populated with a copy of row 285's
# current content (this shifts rows 285:375 down to 286:376, extending the
# sheet's used range to R376).
$ws.Rows("285:285").Copy()
$ws.Rows("285:285").Insert()

# Overwrite the new row 285 with this week's new price observation.
$ws.Cells.Item(285, 4).Value = 44524   # Fecha
$ws.Cells.Item(285, 10).Value = 510    # Volumen
$ws.Cells.Item(285, 12).Value = 9500   # Precio maximo
$ws.Cells.Item(285, 13).Value = 9245   # Precio promedio ponderado
$ws.Cells.Item(285, 16).Value = 370    # Precio $/Kg
